# Auto-generated update script for resum_diari_meteocat.xlsx
# Commit: Update automàtic: dades i banners [2026-02-07 20:19]

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-07 20:17:44'
$ws.Range('E3').Value = '2026-02-07 20:17:47'
$ws.Range('O3').Value = '-5.0 °C'
$ws.Range('E4').Value = '2026-02-07 20:17:49'
$ws.Range('E5').Value = '2026-02-07 20:17:52'
$ws.Range('O5').Value = '10.7 °C'
$ws.Range('E6').Value = '2026-02-07 20:17:55'
$ws.Range('H6').Value = '49%'
$ws.Range('E7').Value = '2026-02-07 20:17:57'
$ws.Range('E8').Value = '2026-02-07 20:17:59'
$ws.Range('H8').Value = '76%'
$ws.Range('O8').Value = '8.6 °C'
$ws.Range('E9').Value = '2026-02-07 20:18:02'
$ws.Range('O9').Value = '3.6 °C'
$ws.Range('E10').Value = '2026-02-07 20:18:04'
$ws.Range('E11').Value = '2026-02-07 20:18:07'
$ws.Range('J11').Value = '1006.4 hPa'
$ws.Range('E12').Value = '2026-02-07 20:18:09'
$ws.Range('H12').Value = '58%'
$ws.Range('O12').Value = '12.1 °C'
$ws.Range('E13').Value = '2026-02-07 20:18:11'
$ws.Range('N13').Value = '5.5 °C 19:59 TU'
$ws.Range('O13').Value = '10.9 °C'
$ws.Range('E14').Value = '2026-02-07 20:18:14'
$ws.Range('O14').Value = '-5.6 °C'
$ws.Range('E15').Value = '2026-02-07 20:18:16'
$ws.Range('E16').Value = '2026-02-07 20:18:19'
$ws.Range('H16').Value = '84%'
$ws.Range('K16').Value = '6.7 MJ/m2'
$ws.Range('E17').Value = '2026-02-07 20:18:21'
$ws.Range('J17').Value = '1006.0 hPa'
$ws.Range('E18').Value = '2026-02-07 20:18:23'
$ws.Range('O18').Value = '-6.1 °C'
$ws.Range('E19').Value = '2026-02-07 20:18:26'
$ws.Range('J19').Value = '1007.3 hPa'
$ws.Range('E20').Value = '2026-02-07 20:18:28'
$ws.Range('E21').Value = '2026-02-07 20:18:31'
$ws.Range('E22').Value = '2026-02-07 20:18:33'
$ws.Range('E23').Value = '2026-02-07 20:18:36'
$ws.Range('H23').Value = '79%'
$ws.Range('O23').Value = '10.1 °C'
$ws.Range('E24').Value = '2026-02-07 20:18:38'
$ws.Range('J24').Value = '1003.3 hPa'
$ws.Range('E25').Value = '2026-02-07 20:18:41'
$ws.Range('J25').Value = '1006.1 hPa'
$ws.Range('O25').Value = '2.3 °C'
$ws.Range('E26').Value = '2026-02-07 20:18:43'
$ws.Range('H26').Value = '70%'
$ws.Range('O26').Value = '-2.3 °C'
$ws.Range('E27').Value = '2026-02-07 20:18:46'
$ws.Range('E28').Value = '2026-02-07 20:18:48'
$ws.Range('E29').Value = '2026-02-07 20:18:51'
$ws.Range('O29').Value = '12.2 °C'
$ws.Range('E30').Value = '2026-02-07 20:18:53'
$ws.Range('O30').Value = '-4.5 °C'
$ws.Range('E31').Value = '2026-02-07 20:18:56'
$ws.Range('J31').Value = '1006.8 hPa'
$ws.Range('E32').Value = '2026-02-07 20:18:58'
$ws.Range('O32').Value = '12.7 °C'
$ws.Range('E33').Value = '2026-02-07 20:19:00'
$ws.Range('O33').Value = '10.4 °C'
$ws.Range('E34').Value = '2026-02-07 20:19:03'
$ws.Range('E35').Value = '2026-02-07 20:19:05'
$ws.Range('E36').Value = '2026-02-07 20:19:08'
$ws.Range('J36').Value = '1007.7 hPa'
